$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue "D2" '26.803.86'
Set-TextValue "E2" '  +1.10%  '
Set-TextValue "D3" '1.649.58'
Set-TextValue "E3" '  +1.37%  '
Set-TextValue "E4" '  +0.75%  '
Set-TextValue "D5" '216.57'
Set-TextValue "E5" '  +1.65%  '
Set-TextValue "E6" '  +0.89%  '
Set-TextValue "E7" '  +0.63%  '
Set-TextValue "D9" '0.0627'
Set-TextValue "E9" '  +0.70%  '
Set-TextValue "E11" '  +0.07%  '
Set-TextValue "D12" '1.879.79'
Set-TextValue "E12" '  +1.40%  '
Set-TextValue "D13" '1.645.81'
Set-TextValue "E13" '  +0.89%  '
Set-TextValue "E14" '  +1.48%  '
Set-TextValue "E15" '  +2.01%  '
Set-TextValue "D16" '65.43'
Set-TextValue "E16" '  +0.64%  '
Set-TextValue "D17" '26.799.28'
Set-TextValue "E17" '  +0.96%  '
Set-TextValue "D18" '0.0₃0744'
Set-TextValue "E18" '  +0.60%  '
Set-TextValue "D19" '218.38'
Set-TextValue "E19" '  +1.85%  '
Set-TextValue "E21" '  +1.71%  '
Set-TextValue "D22" '2.44'
Set-TextValue "E22" '  +17.88%  '
Set-TextValue "D23" '6.26'
Set-TextValue "E23" '  +0.16%  '
Set-TextValue "E24" '  +2.27%  '
Set-TextValue "D25" '146.46'
Set-TextValue "E25" '  -1.27%  '
Set-TextValue "E26" '  +0.45%  '
Set-TextValue "E27" '  -0.08%  '
Set-TextValue "E28" '  +3.82%  '
Set-TextValue "D29" '15.74'
Set-TextValue "E29" '  +1.41%  '
Set-TextValue "E30" '  +1.80%  '
Set-TextValue "E31" '  +1.84%  '
Set-TextValue "D32" '3.34'
Set-TextValue "E32" '  +0.36%  '
Set-TextValue "E33" '  +1.60%  '
Set-TextValue "D34" '1.283.87'
Set-TextValue "E34" '  +3.69%  '
Set-TextValue "E35" '  +2.95%  '
Set-TextValue "E36" '  +3.11%  '
Set-TextValue "E37" '  +2.62%  '
Set-TextValue "D38" '0.537'
Set-TextValue "E38" '  +5.89%  '
Set-TextValue "E39" '  +4.09%  '
Set-TextValue "E40" '  +0.61%  '
Set-TextValue "D41" '0.814'
Set-TextValue "E41" '  +1.85%  '
Set-TextValue "E42" '  -0.83%  '
Set-TextValue "E43" '  +2.56%  '
Set-TextValue "D44" '1.789.58'
Set-TextValue "E44" '  +1.40%  '
Set-TextValue "D45" '91.90'
Set-TextValue "E45" '  -1.19%  '
Set-TextValue "D46" '59.65'
Set-TextValue "E46" '  +8.73%  '
Set-TextValue "E47" '  +1.54%  '
Set-TextValue "E48" '  +1.38%  '
Set-TextValue "D49" '7.77'
Set-TextValue "E49" '  +3.93%  '
Set-TextValue "D50" '0.0969'
Set-TextValue "E50" '  +1.84%  '
Set-TextValue "E51" '  +0.68%  '
